$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its current location (end of the Due
#    Date line, right after the "Dec" run).
# ---------------------------------------------------------------------------
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
}

# ---------------------------------------------------------------------------
# 2. Locate the "References" hyperlink paragraph
#      (<w:hyperlink r:id="rId5">...</w:hyperlink><w:hyperlink r:id="rId6"/>)
#    by searching for its visible text, rather than assuming a fixed
#    paragraph index.
# ---------------------------------------------------------------------------

$findRng = $d.Content
$found = $findRng.Find.Execute("https://azure.microsoft.com/en-us/get-started/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$pLinkIndex = 0
$idx = 0
foreach ($pp in $d.Paragraphs) {
    $idx = $idx + 1
    if ($pp.Range.Start -le $findRng.Start -and $pp.Range.End -ge $findRng.End) {
        $pLinkIndex = $idx
    }
}

$pLink = $d.Paragraphs($pLinkIndex)
$pEmpty = $d.Paragraphs($pLinkIndex + 1)

# ---------------------------------------------------------------------------
# 3. Restructure those two paragraphs into three:
#      a) new paragraph, new pPr, hyperlink -> data.zip (rId5, reused)
#      b) paragraph with jc=both containing the still-empty rId6 hyperlink
#      c) empty paragraph now hosting the _GoBack bookmark
#
#    Done in two InsertXML passes because a single InsertXML call that
#    mints 2+ brand-new paragraphs at once drops <w:rStyle> from the
#    inserted runs; splitting the structural insert from the rStyle-
#    bearing run avoids that.
# ---------------------------------------------------------------------------

$rStart = $pLink.Range.Start
$rEnd = $pEmpty.Range.End
$rng = $d.Range($rStart, $rEnd)

$structureXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
<w:pPr>
<w:tabs><w:tab w:val="clear" w:pos="720"/></w:tabs>
<w:suppressAutoHyphens w:val="0"/>
<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr>
</w:pPr>
<w:hyperlink r:id="rId5" w:history="1"><w:r><w:t>https://xuhappy.github.io/courses/BigData/homework/data.zip</w:t></w:r></w:hyperlink>
</w:p>
<w:p>
<w:pPr><w:jc w:val="both"/></w:pPr>
<w:hyperlink r:id="rId6"></w:hyperlink>
</w:p>
<w:p>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng.InsertXML($structureXml)

# ---------------------------------------------------------------------------
# 4. Now that the three paragraphs exist, replace the first one's content
#     again (single-paragraph InsertXML) so the run picks up rStyle=Hyperlink
#     without losing it to the multi-paragraph-insert quirk.
# ---------------------------------------------------------------------------

$pNewLink = $d.Paragraphs($pLinkIndex)
$rng2 = $d.Range($pNewLink.Range.Start, $pNewLink.Range.End)

$styleXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
<w:pPr>
<w:tabs><w:tab w:val="clear" w:pos="720"/></w:tabs>
<w:suppressAutoHyphens w:val="0"/>
<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr>
</w:pPr>
<w:hyperlink r:id="rId5" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://xuhappy.github.io/courses/BigData/homework/data.zip</w:t></w:r></w:hyperlink>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng2.InsertXML($styleXml)

# ---------------------------------------------------------------------------
# 5. Point the reused hyperlink relationship at the new target URL.
# ---------------------------------------------------------------------------

$pFinalLink = $d.Paragraphs($pLinkIndex)
$h = $pFinalLink.Range.Hyperlinks.Item(1)
$h.Address = "https://xuhappy.github.io/courses/BigData/homework/data.zip"

Write-Host "Done. Paragraphs:" $d.Paragraphs.Count
